$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 phone number was stored as text; normalize it to a real number
$ws.Range("A23").Value = 79174445

# Append the new redemption row: "Redeem points 71076783 100.0"
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "71076783"
$ws.Range("B24").Value = 100
$ws.Range("C24").Value = "2025-08-18T17:50:48"
